# Auto-generated Excel COM-interop edit script
# Updates Leve market-price / profit figures across 8 worksheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match refreshed API data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 448.25
$ws.Range("I39").Value = 341
$ws.Range("J39").Value = 1199
$ws.Range("K39").Value = 1023
$ws.Range("L39").Value = 3597
$ws.Range("M39").Value = -727
$ws.Range("N39").Value = -4189

$ws.Range("H132").Value = 1526.44
$ws.Range("I132").Value = 1527.5416
$ws.Range("K132").Value = 4582.6248
$ws.Range("M132").Value = -2052.6248

$ws.Range("H138").Value = 6151.741
$ws.Range("I138").Value = 7451.905
$ws.Range("J138").Value = 5413.811
$ws.Range("K138").Value = 22355.715
$ws.Range("L138").Value = 16241.433
$ws.Range("M138").Value = -17215.715
$ws.Range("N138").Value = -26521.433

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20158.38
$ws.Range("I32").Value = 19956.646
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 19956.646
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -19669.646
$ws.Range("N32").Value = -25574

$ws.Range("H45").Value = 3698.75
$ws.Range("I45").Value = 2754.8125
$ws.Range("K45").Value = 2754.8125
$ws.Range("M45").Value = -2377.8125

$ws.Range("H63").Value = 4803
$ws.Range("I63").Value = 4637
$ws.Range("J63").Value = 4969
$ws.Range("K63").Value = 4637
$ws.Range("L63").Value = 4969
$ws.Range("M63").Value = -3951
$ws.Range("N63").Value = -6341

$ws.Range("H66").Value = 4803
$ws.Range("I66").Value = 4637
$ws.Range("J66").Value = 4969
$ws.Range("K66").Value = 23185
$ws.Range("L66").Value = 24845
$ws.Range("M66").Value = -19753
$ws.Range("N66").Value = -31709

$ws.Range("H74").Value = 1884.8334
$ws.Range("I74").Value = 931.2632
$ws.Range("J74").Value = 5508.4
$ws.Range("K74").Value = 931.2632
$ws.Range("L74").Value = 5508.4
$ws.Range("M74").Value = -57.26319999999998
$ws.Range("N74").Value = -7256.4

$ws.Range("H77").Value = 1884.8334
$ws.Range("I77").Value = 931.2632
$ws.Range("J77").Value = 5508.4
$ws.Range("K77").Value = 4656.316
$ws.Range("L77").Value = 27542
$ws.Range("M77").Value = -288.3159999999998
$ws.Range("N77").Value = -36278

$ws.Range("H132").Value = 2190.6167
$ws.Range("I132").Value = 2168.868
$ws.Range("K132").Value = 6506.603999999999
$ws.Range("M132").Value = -3976.603999999999

$ws.Range("H135").Value = 82206.836
$ws.Range("J135").Value = 82206.836
$ws.Range("L135").Value = 82206.836
$ws.Range("N135").Value = -92346.836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2184.5
$ws.Range("I99").Value = 999
$ws.Range("K99").Value = 999
$ws.Range("M99").Value = 499

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H140").Value = 83984.5
$ws.Range("I140").Value = 83984.5
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 83984.5
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -78804.5
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5951.6665
$ws.Range("I31").Value = 4499.25
$ws.Range("J31").Value = 6479.8184
$ws.Range("K31").Value = 4499.25
$ws.Range("L31").Value = 6479.8184
$ws.Range("M31").Value = -4204.25
$ws.Range("N31").Value = -7069.8184

$ws.Range("H34").Value = 5951.6665
$ws.Range("I34").Value = 4499.25
$ws.Range("J34").Value = 6479.8184
$ws.Range("K34").Value = 4499.25
$ws.Range("L34").Value = 6479.8184
$ws.Range("M34").Value = -4297.25
$ws.Range("N34").Value = -6883.8184

$ws.Range("H58").Value = 5127.857
$ws.Range("I58").Value = 2917.875
$ws.Range("J58").Value = 12199.8
$ws.Range("K58").Value = 2917.875
$ws.Range("L58").Value = 12199.8
$ws.Range("M58").Value = -2714.875
$ws.Range("N58").Value = -12605.8

$ws.Range("H62").Value = 6681.143
$ws.Range("I62").Value = 5474.75
$ws.Range("K62").Value = 5474.75
$ws.Range("M62").Value = -4850.75

$ws.Range("H65").Value = 6681.143
$ws.Range("I65").Value = 5474.75
$ws.Range("K65").Value = 27373.75
$ws.Range("M65").Value = -24253.75

$ws.Range("H82").Value = 54999
$ws.Range("J82").Value = 54999
$ws.Range("L82").Value = 54999
$ws.Range("N82").Value = -55721

$ws.Range("H85").Value = 54999
$ws.Range("J85").Value = 54999
$ws.Range("L85").Value = 54999
$ws.Range("N85").Value = -57495

$ws.Range("H99").Value = 7945.846
$ws.Range("J99").Value = 8999.625
$ws.Range("L99").Value = 8999.625
$ws.Range("N99").Value = -11995.625

$ws.Range("H122").Value = 2382.1292
$ws.Range("I122").Value = 2417.45
$ws.Range("J122").Value = 2317.9092
$ws.Range("K122").Value = 7252.349999999999
$ws.Range("L122").Value = 6953.7276
$ws.Range("M122").Value = -4802.349999999999
$ws.Range("N122").Value = -11853.7276

$ws.Range("H126").Value = 7945.846
$ws.Range("J126").Value = 8999.625
$ws.Range("L126").Value = 26998.875
$ws.Range("N126").Value = -31938.875

$ws.Range("H134").Value = 2590.9744
$ws.Range("I134").Value = 1531.3334
$ws.Range("K134").Value = 4594.0002
$ws.Range("M134").Value = -2059.0002

$ws.Range("H136").Value = 5127.857
$ws.Range("I136").Value = 2917.875
$ws.Range("J136").Value = 12199.8
$ws.Range("K136").Value = 8753.625
$ws.Range("L136").Value = 36599.39999999999
$ws.Range("M136").Value = -6203.625
$ws.Range("N136").Value = -41699.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1085.6
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1085.6
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3256.8
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -3602.8

$ws.Range("H114").Value = 2190
$ws.Range("I114").Value = 1253.3334
$ws.Range("J114").Value = 5000
$ws.Range("K114").Value = 3760.0002
$ws.Range("L114").Value = 15000
$ws.Range("M114").Value = -506.0001999999999
$ws.Range("N114").Value = -21508

$ws.Range("H122").Value = 11810.889
$ws.Range("J122").Value = 12672.565
$ws.Range("L122").Value = 114053.085
$ws.Range("N122").Value = -118953.085

$ws.Range("H129").Value = 339837.56
$ws.Range("J129").Value = 599704.2
$ws.Range("L129").Value = 1799112.6
$ws.Range("N129").Value = -1809112.6

$ws.Range("H131").Value = 3550.3462
$ws.Range("I131").Value = 1808.875
$ws.Range("K131").Value = 5426.625
$ws.Range("M131").Value = -386.625

$ws.Range("H132").Value = 2359.375
$ws.Range("J132").Value = 2359.375
$ws.Range("L132").Value = 21234.375
$ws.Range("N132").Value = -26294.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3355.875
$ws.Range("I31").Value = 3355.875
$ws.Range("K31").Value = 3355.875
$ws.Range("M31").Value = -3063.875

$ws.Range("H37").Value = 3355.875
$ws.Range("I37").Value = 3355.875
$ws.Range("K37").Value = 3355.875
$ws.Range("M37").Value = -3078.875

$ws.Range("H122").Value = 2851.75
$ws.Range("J122").Value = 1833.2222
$ws.Range("L122").Value = 5499.6666
$ws.Range("N122").Value = -10399.6666

$ws.Range("H135").Value = 147999.67
$ws.Range("J135").Value = 147999.67
$ws.Range("L135").Value = 147999.67
$ws.Range("N135").Value = -158139.67

$ws.Range("H137").Value = 63558.168
$ws.Range("I137").Value = 63558.168
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 63558.168
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -58458.168
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 7559.5713
$ws.Range("I3").Value = 7559.5713
$ws.Range("K3").Value = 7559.5713
$ws.Range("M3").Value = -7447.5713

$ws.Range("H15").Value = 7559.5713
$ws.Range("I15").Value = 7559.5713
$ws.Range("K15").Value = 7559.5713
$ws.Range("M15").Value = -7389.5713

$ws.Range("H22").Value = 7099.2
$ws.Range("J22").Value = 8928.700000000001
$ws.Range("L22").Value = 8928.700000000001
$ws.Range("N22").Value = -9518.700000000001

$ws.Range("H27").Value = 7099.2
$ws.Range("J27").Value = 8928.700000000001
$ws.Range("L27").Value = 8928.700000000001
$ws.Range("N27").Value = -9142.700000000001

$ws.Range("H68").Value = 6579.8
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251

$ws.Range("H71").Value = 6579.8
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256

$ws.Range("H122").Value = 5722
$ws.Range("I122").Value = 6366.4
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 19099.2
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -16649.2
$ws.Range("N122").Value = -12400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

$ws.Range("H136").Value = 3880.5312
$ws.Range("I136").Value = 2228.8262
$ws.Range("J136").Value = 8101.5557
$ws.Range("K136").Value = 6686.4786
$ws.Range("L136").Value = 24304.6671
$ws.Range("M136").Value = -4136.4786
$ws.Range("N136").Value = -29404.6671

$ws.Range("H141").Value = 80701.664
$ws.Range("J141").Value = 80701.664
$ws.Range("L141").Value = 80701.664
$ws.Range("N141").Value = -91061.664
